# Refresh monthly data-vintage figures (M2/FX series lengths & date ranges)
# for Top50_DataComp as of the new data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 30074
$ws.Range("H2").Value = 45261

# Row 3
$ws.Range("E3").Value = 30011
$ws.Range("F3").Value = 45200

# Row 4
$ws.Range("E4").Value = 30011
$ws.Range("F4").Value = 45200
$ws.Range("G4").Value = 30074
$ws.Range("H4").Value = 45261

# Row 5
$ws.Range("G5").Value = 30074
$ws.Range("H5").Value = 45261

# Row 6
$ws.Range("C6").Value = 443
$ws.Range("F6").Value = 45200
$ws.Range("G6").Value = 30074
$ws.Range("H6").Value = 45261

# Row 7
$ws.Range("G7").Value = 30074
$ws.Range("H7").Value = 45261

# Row 8
$ws.Range("D8").Value = 410
$ws.Range("E8").Value = 30011
$ws.Range("F8").Value = 45200
$ws.Range("H8").Value = 45261

# Row 9
$ws.Range("E9").Value = 30011
$ws.Range("F9").Value = 45200
$ws.Range("G9").Value = 30074
$ws.Range("H9").Value = 45261

# Row 10
$ws.Range("D10").Value = 483
$ws.Range("E10").Value = 30011
$ws.Range("F10").Value = 45200
$ws.Range("H10").Value = 45261

# Row 11
$ws.Range("G11").Value = 30074
$ws.Range("H11").Value = 45261

# Row 12
$ws.Range("C12").Value = 371
$ws.Range("D12").Value = 352
$ws.Range("F12").Value = 45200
$ws.Range("H12").Value = 45261

# Row 13
$ws.Range("C13").Value = 467
$ws.Range("F13").Value = 45200
$ws.Range("G13").Value = 30074
$ws.Range("H13").Value = 45261

# Row 14
$ws.Range("D14").Value = 396
$ws.Range("H14").Value = 45261

# Row 15
$ws.Range("G15").Value = 30042
$ws.Range("H15").Value = 45261

# Row 16
$ws.Range("C16").Value = 455
$ws.Range("D16").Value = 410
$ws.Range("F16").Value = 45200
$ws.Range("H16").Value = 45261

# Row 17
$ws.Range("C17").Value = 370
$ws.Range("D17").Value = 394
$ws.Range("F17").Value = 45200
$ws.Range("H17").Value = 45261

# Row 18
$ws.Range("E18").Value = 30011
$ws.Range("F18").Value = 45200
$ws.Range("G18").Value = 30074
$ws.Range("H18").Value = 45261

# Row 19
$ws.Range("D19").Value = 398
$ws.Range("E19").Value = 30011
$ws.Range("F19").Value = 45200
$ws.Range("H19").Value = 45261

# Row 20
$ws.Range("C20").Value = 481
$ws.Range("F20").Value = 45200
$ws.Range("G20").Value = 30074
$ws.Range("H20").Value = 45261

# Row 21
$ws.Range("C21").Value = 310
$ws.Range("F21").Value = 45200
$ws.Range("G21").Value = 30074
$ws.Range("H21").Value = 45261

# Row 22
$ws.Range("C22").Value = 323
$ws.Range("D22").Value = 367
$ws.Range("F22").Value = 45200
$ws.Range("H22").Value = 45261

# Row 23
$ws.Range("D23").Value = 402
$ws.Range("H23").Value = 45261

# Row 24
$ws.Range("C24").Value = 216
$ws.Range("D24").Value = 410
$ws.Range("F24").Value = 45231
$ws.Range("H24").Value = 45261

# Row 25
$ws.Range("D25").Value = 315
$ws.Range("H25").Value = 45261

# Row 26
$ws.Range("C26").Value = 333
$ws.Range("D26").Value = 313
$ws.Range("F26").Value = 45200
$ws.Range("H26").Value = 45261

# Row 27
$ws.Range("E27").Value = 30011
$ws.Range("F27").Value = 45200
$ws.Range("G27").Value = 30074
$ws.Range("H27").Value = 45261

# Row 28
$ws.Range("D28").Value = 380
$ws.Range("H28").Value = 45261

# Row 29
$ws.Range("C29").Value = 262
$ws.Range("D29").Value = 235
$ws.Range("F29").Value = 45200
$ws.Range("H29").Value = 45261

# Row 30
$ws.Range("D30").Value = 217
$ws.Range("E30").Value = 30011
$ws.Range("F30").Value = 45200
$ws.Range("H30").Value = 45261

# Row 31
$ws.Range("C31").Value = 394
$ws.Range("F31").Value = 45200
$ws.Range("G31").Value = 30074
$ws.Range("H31").Value = 45261

# Row 32
$ws.Range("E32").Value = 30011
$ws.Range("F32").Value = 45200
$ws.Range("G32").Value = 30074
$ws.Range("H32").Value = 45261

# Row 33
$ws.Range("D33").Value = 398
$ws.Range("H33").Value = 45261

# Row 34
$ws.Range("C34").Value = 202
$ws.Range("D34").Value = 317
$ws.Range("F34").Value = 45200
$ws.Range("H34").Value = 45261

# Row 35
$ws.Range("C35").Value = 404
$ws.Range("D35").Value = 317
$ws.Range("F35").Value = 45170
$ws.Range("H35").Value = 45261

# Row 36
$ws.Range("D36").Value = 410
$ws.Range("H36").Value = 45261

# Row 37
$ws.Range("C37").Value = 466
$ws.Range("D37").Value = 317
$ws.Range("F37").Value = 45200
$ws.Range("H37").Value = 45261

# Row 38
$ws.Range("C38").Value = 359
$ws.Range("D38").Value = 367
$ws.Range("F38").Value = 45200
$ws.Range("H38").Value = 45261

# Row 39
$ws.Range("C39").Value = 227
$ws.Range("D39").Value = 223
$ws.Range("F39").Value = 45200
$ws.Range("H39").Value = 45261

# Row 40
$ws.Range("D40").Value = 314
$ws.Range("H40").Value = 45261

# Row 41
$ws.Range("D41").Value = 316
$ws.Range("H41").Value = 45261

# Row 42
$ws.Range("D42").Value = 218
$ws.Range("H42").Value = 45261

# Row 43
$ws.Range("D43").Value = 317
$ws.Range("H43").Value = 45261

# Row 44
$ws.Range("D44").Value = 304
$ws.Range("H44").Value = 45261

# Row 45
$ws.Range("D45").Value = 317
$ws.Range("H45").Value = 45261

# Row 46
$ws.Range("C46").Value = 334
$ws.Range("D46").Value = 298
$ws.Range("F46").Value = 45200
$ws.Range("H46").Value = 45261

# Row 47
$ws.Range("C47").Value = 335
$ws.Range("D47").Value = 260
$ws.Range("F47").Value = 45200
$ws.Range("H47").Value = 45261

# Row 48
$ws.Range("D48").Value = 315
$ws.Range("H48").Value = 45261

# Row 49
$ws.Range("D49").Value = 313
$ws.Range("H49").Value = 45261

# Row 50
$ws.Range("D50").Value = 238
$ws.Range("H50").Value = 45261

# Row 51
$ws.Range("D51").Value = 317
$ws.Range("E51").Value = 29830
$ws.Range("F51").Value = 45170
$ws.Range("H51").Value = 45261

# Row 52
$ws.Range("D52").Value = 315
$ws.Range("H52").Value = 45261

